# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets.
$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

# Map of row number -> new value for column F
$updates = @{
    2  = 11670
    3  = 11264
    4  = 605
    11 = 10725
    12 = 4148
    14 = 5
    16 = 2464
    19 = 124
    20 = 444
    21 = 11131
    22 = 10909
}

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
